# Update the dSF column (F) with repulled / recalculated values for several rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 5
    3  = 1
    4  = 5
    6  = 2
    10 = -1
    12 = -8
    13 = -2
    14 = -7
    15 = -7
    17 = 3
    24 = -10
    25 = -6
    27 = 4
    29 = -4
    32 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
